$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 32 (G=5484)
$ws.Range("H32").Value = 605.8125
$ws.Range("I32").Value = 660
$ws.Range("J32").Value = 587.75
$ws.Range("K32").Value = 660
$ws.Range("L32").Value = 587.75
$ws.Range("M32").Value = -334
$ws.Range("N32").Value = -1239.75

# row 41 (G=5478)
$ws.Range("H41").Value = 82.14286
$ws.Range("I41").Value = 83.75
$ws.Range("J41").Value = 80
$ws.Range("K41").Value = 83.75
$ws.Range("L41").Value = 80
$ws.Range("M41").Value = 356.25
$ws.Range("N41").Value = -960

# row 74 (G=5507)
$ws.Range("H74").Value = 3453.0417
$ws.Range("I74").Value = 2945
$ws.Range("K74").Value = 2945
$ws.Range("M74").Value = -2009

# row 76 (G=12602)
$ws.Range("H76").Value = 6721.478
$ws.Range("I76").Value = 4940.8
$ws.Range("J76").Value = 7216.1113
$ws.Range("K76").Value = 4940.8
$ws.Range("L76").Value = 7216.1113
$ws.Range("M76").Value = -4625.8
$ws.Range("N76").Value = -7846.1113

# row 77 (G=5507)
$ws.Range("H77").Value = 3453.0417
$ws.Range("I77").Value = 2945
$ws.Range("K77").Value = 14725
$ws.Range("M77").Value = -10045

# row 79 (G=12602)
$ws.Range("H79").Value = 6721.478
$ws.Range("I79").Value = 4940.8
$ws.Range("J79").Value = 7216.1113
$ws.Range("K79").Value = 4940.8
$ws.Range("L79").Value = 7216.1113
$ws.Range("M79").Value = -3848.8
$ws.Range("N79").Value = -9400.1113

# row 137 (G=44013)
$ws.Range("H137").Value = 4607.6665
$ws.Range("I137").Value = 3124.8333
$ws.Range("J137").Value = 4854.8057
$ws.Range("K137").Value = 9374.499899999999
$ws.Range("L137").Value = 14564.4171
$ws.Range("M137").Value = -6824.499899999999
$ws.Range("N137").Value = -19664.4171

# row 138 (G=44169)
$ws.Range("H138").Value = 3672.8386
$ws.Range("I138").Value = 2425.5715
$ws.Range("J138").Value = 4700
$ws.Range("K138").Value = 7276.7145
$ws.Range("L138").Value = 14100
$ws.Range("M138").Value = -2136.7145
$ws.Range("N138").Value = -24380

$ws = $wb.Worksheets.Item("ARM")
# row 32 (G=44147)
$ws.Range("H32").Value = 13787.826
$ws.Range("I32").Value = 11300.825
$ws.Range("J32").Value = 30367.834
$ws.Range("K32").Value = 11300.825
$ws.Range("L32").Value = 30367.834
$ws.Range("M32").Value = -11013.825
$ws.Range("N32").Value = -30941.834

# row 63 (G=12528)
$ws.Range("H63").Value = 1928.3334
$ws.Range("I63").Value = 1928.3334
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1928.3334
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1242.3334
$ws.Range("N63").ClearContents()

# row 66 (G=12528)
$ws.Range("H66").Value = 1928.3334
$ws.Range("I66").Value = 1928.3334
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9641.666999999999
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6209.666999999999
$ws.Range("N66").ClearContents()

# row 92 (G=18050)
$ws.Range("H92").Value = 37800
$ws.Range("J92").Value = 37800
$ws.Range("L92").Value = 37800
$ws.Range("N92").Value = -42792

$ws = $wb.Worksheets.Item("CRP")
# row 92 (G=18041)
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# row 103 (G=19558)
$ws.Range("H103").Value = 5770.6665
$ws.Range("I103").Value = 5770.6665
$ws.Range("K103").Value = 5770.6665
$ws.Range("M103").Value = -4598.6665

$ws = $wb.Worksheets.Item("CUL")
# row 23 (G=4858)
$ws.Range("H23").Value = 80.09999999999999
$ws.Range("I23").Value = 102.75
$ws.Range("J23").Value = 65
$ws.Range("K23").Value = 308.25
$ws.Range("L23").Value = 195
$ws.Range("M23").Value = -73.25
$ws.Range("N23").Value = -665

# row 32 (G=4731)
$ws.Range("H32").Value = 290440
$ws.Range("I32").Value = 322488.88
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 967466.64
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -967183.64
$ws.Range("N32").Value = -6566

# row 68 (G=12895)
$ws.Range("H68").Value = 841.26666
$ws.Range("I68").Value = 835.6429000000001
$ws.Range("J68").Value = 846.1875
$ws.Range("K68").Value = 2506.9287
$ws.Range("L68").Value = 2538.5625
$ws.Range("M68").Value = -1695.9287
$ws.Range("N68").Value = -4160.5625

# row 71 (G=12895)
$ws.Range("H71").Value = 841.26666
$ws.Range("I71").Value = 835.6429000000001
$ws.Range("J71").Value = 846.1875
$ws.Range("K71").Value = 7520.7861
$ws.Range("L71").Value = 7615.6875
$ws.Range("M71").Value = -3464.7861
$ws.Range("N71").Value = -15727.6875

# row 80 (G=12890)
$ws.Range("H80").Value = 86483.25
$ws.Range("I80").Value = 169333.17
$ws.Range("J80").Value = 3633.3333
$ws.Range("K80").Value = 507999.51
$ws.Range("L80").Value = 10899.9999
$ws.Range("M80").Value = -507063.51
$ws.Range("N80").Value = -12771.9999

# row 83 (G=12890)
$ws.Range("H83").Value = 86483.25
$ws.Range("I83").Value = 169333.17
$ws.Range("J83").Value = 3633.3333
$ws.Range("K83").Value = 1523998.53
$ws.Range("L83").Value = 32699.9997
$ws.Range("M83").Value = -1519318.53
$ws.Range("N83").Value = -42059.9997

# row 139 (G=44102)
$ws.Range("H139").Value = 251841.02
$ws.Range("I139").Value = 478943.88
$ws.Range("J139").Value = 3109.3333
$ws.Range("K139").Value = 1436831.64
$ws.Range("L139").Value = 9327.999899999999
$ws.Range("M139").Value = -1431691.64
$ws.Range("N139").Value = -19607.9999

$ws = $wb.Worksheets.Item("GSM")
# row 80 (G=12521)
$ws.Range("H80").Value = 3300
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3300
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3300
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -5296

# row 83 (G=12521)
$ws.Range("H83").Value = 3300
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3300
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 16500
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -26484

$ws = $wb.Worksheets.Item("LTW")
# row 7 (G=36249)
$ws.Range("H7").Value = 2094.25
$ws.Range("I7").Value = 1386
$ws.Range("J7").Value = 2802.5
$ws.Range("K7").Value = 1386
$ws.Range("L7").Value = 2802.5
$ws.Range("M7").Value = -1274
$ws.Range("N7").Value = -3026.5

# row 126 (G=36249)
$ws.Range("H126").Value = 2094.25
$ws.Range("I126").Value = 1386
$ws.Range("J126").Value = 2802.5
$ws.Range("K126").Value = 4158
$ws.Range("L126").Value = 8407.5
$ws.Range("M126").Value = -1688
$ws.Range("N126").Value = -13347.5

$ws = $wb.Worksheets.Item("WVR")
# row 70 (G=11979)
$ws.Range("H70").Value = 36250
$ws.Range("I70").Value = 20000
$ws.Range("J70").Value = 41666.668
$ws.Range("K70").Value = 20000
$ws.Range("L70").Value = 41666.668
$ws.Range("M70").Value = -19685
$ws.Range("N70").Value = -42296.668

# row 73 (G=11979)
$ws.Range("H73").Value = 36250
$ws.Range("I73").Value = 20000
$ws.Range("J73").Value = 41666.668
$ws.Range("K73").Value = 20000
$ws.Range("L73").Value = 41666.668
$ws.Range("M73").Value = -20000
$ws.Range("N73").Value = -43684
